# otorgantes.docx edit:
#  - "+++HTML documento.minuta+++" template tag -> "+++= documento.minuta+++",
#    now rendered bold
#  - static "minuta que está firmada por el Abogado FERNANDO BARRERA. MAT.
#    18-2015-199," text removed
#  - the fixed phrase "le fue a la compareciente" (right after "leída que ")
#    is replaced by a conditional template expression choosing between
#    singular/plural and gendered wording
#  - the _GoBack bookmark, which used to wrap "+++HTML documento.minuta+++",
#    now wraps the (empty) point right after the new template expression

$d = $word.ActiveDocument

# 1) "+++HTML " -> "+++= " (the two runs "+++HTML" and " " merge into one
#    run with text "+++= ")
$d.Content.Find.Execute("+++HTML ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "+++= ", 2)

# 2) Make the whole "+++= documento.minuta+++" template tag bold
$r = $d.Content
$r.Find.Execute("+++= documento.minuta+++", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$r.Bold = 1

# 3) The _GoBack bookmark used to sit around that template tag; it is going
#    to be re-created further down, right after the new template expression
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 4) Drop "minuta que está firmada por el Abogado FERNANDO BARRERA. MAT.
#    18-2015-199, " (including the bold lawyer name run) entirely
$d.Content.Find.Execute("minuta que está firmada por el Abogado FERNANDO BARRERA. MAT. 18-2015-199, ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 5) Replace the fixed "le fue a la compareciente " wording with the new
#    conditional template expression (note: use single-quoted PowerShell
#    strings throughout so the back-ticks in the template text are taken
#    literally instead of being treated as escape characters)
$newMiddle = '+++= documento.otorgantes.length > 1? ` les fue a los comparecientes `: documento.otorgantes[0].tratamiento == `EL SEÑOR`? ` le fue al compareciente ` : ` le fue a la compareciente`+++ '
$d.Content.Find.Execute('le fue a la compareciente ', $false, $false, $false, `
    $false, $false, $true, 1, $false, $newMiddle, 2)

# 6) Re-create the _GoBack bookmark as a zero-length bookmark right after
#    the closing "+++" of the new template expression (i.e. right before
#    the trailing " por  mí  el  notario..." text)
$r2 = $d.Content
$r2.Find.Execute('le fue a la compareciente`+++', $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)
